$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Insert 5 new worker rows below the existing one (row 16), pushing the
# signature block (old rows 21-22) down to rows 26-27 ---
$ws.Rows("17:21").Insert()

# Clone the formatting of the existing data row (16) onto the five new rows
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J21").PasteSpecial()
$excel.CutCopyMode = $false

# --- Update header summary numbers ---
$ws.Range("E11").Value = 294190
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 2

# --- Fill in the new worker detail rows ---
# Row 17: YOSNEIDER DE JESUS TEHERAN CASTRO
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1002320841"
$ws.Range("D17").Value = "YOSNEIDER DE JESUS TEHERAN CASTRO"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18: RAFAEL EMILIO MARQUEZ GONZALEZ
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "10951468"
$ws.Range("D18").Value = "RAFAEL EMILIO MARQUEZ GONZALEZ"
$ws.Range("E18").Value = "2509"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19: CAROLINA MONTOYA ZULETA
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1080936207"
$ws.Range("D19").Value = "CAROLINA MONTOYA ZULETA"
$ws.Range("E19").Value = "2509"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# Row 20: JUAN JOSE MARQUEZ MOLINA
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1104407196"
$ws.Range("D20").Value = "JUAN JOSE MARQUEZ MOLINA"
$ws.Range("E20").Value = "2509"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

# Row 21: JORGE ERMIRO MARQUEZ MOLINA (repeat worker, new period)
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1100082519"
$ws.Range("D21").Value = "JORGE ERMIRO MARQUEZ MOLINA"
$ws.Range("E21").Value = "2509"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

# --- Column D needs to be widened to fit the new, longer name ---
$ws.Columns("D").AutoFit()
